$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45, shifting existing rows 45:95 down to 46:96
$ws.Rows("45").Insert()

# Populate the new row 45 with the new data entry (copy of row pattern, new values)
$ws.Cells.Item(45, 1).Value = 10
$ws.Cells.Item(45, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(45, 3).Value = "La Araucanía"
$ws.Cells.Item(45, 4).Value = 45079
$ws.Cells.Item(45, 4).NumberFormat = $ws.Cells.Item(46, 4).NumberFormat
$ws.Cells.Item(45, 5).Value = 9
$ws.Cells.Item(45, 6).Value = "Fruta"
$ws.Cells.Item(45, 7).Value = 100108
$ws.Cells.Item(45, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(45, 9).Value = 100108003
$ws.Cells.Item(45, 10).Value = "Maracuyá"
$ws.Cells.Item(45, 11).Value = "Sin especificar"
$ws.Cells.Item(45, 12).Value = "Primera"
$ws.Cells.Item(45, 13).Value = 20
$ws.Cells.Item(45, 14).Value = 50000
$ws.Cells.Item(45, 15).Value = 50000
$ws.Cells.Item(45, 16).Value = 50000
$ws.Cells.Item(45, 17).Value = "$/caja 18 kilos"
$ws.Cells.Item(45, 18).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(45, 19).Value = 2778
$ws.Cells.Item(45, 20).Value = 18
